# Daily attendance processing - 2025-10-13 03:45:02
#
# The "Recorded By" column (G) lists the users who touched an attendance
# session, separated by ", ". Reorder the trailing two entries for every
# row whose text exactly matches one of the two known "before" values:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, system, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#
# Only the rows that actually carry one of those two exact strings are
# touched; every other cell (including the many blank "Recorded By"
# cells for still-pending sessions) is left completely untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 6, 10, 11, 12, 13, 14, 15, 17, 18, 29, 30, 33, 37, 38, 39, 40, 41, 42, 44, 45, 56, 57, 60, 64, 65, 66, 67, 68, 69, 71, 72, 86, 87, 88, 89, 93, 95, 96, 112, 113, 114, 115, 119, 121, 122, 138, 139, 140, 141, 145, 147, 148)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $current = $cell.Value2

    if ($current -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($current -eq "System, system, backup@backdoor.com") {
        $cell.Value = "System, backup@backdoor.com, system"
    }
}
